$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-13 Monday" "2024-05-14 Tuesday"

Replace-Text "491×6=" "121×2="
Replace-Text "134×4=" "829×6="
Replace-Text "740×8=" "703×5="
Replace-Text "686×8=" "985×8="
Replace-Text "952×9=" "952×2="
Replace-Text "550×2=" "392×9="
Replace-Text "305×3=" "416×2="
Replace-Text "273×4=" "438×9="
Replace-Text "617×8=" "708×8="
Replace-Text "750×6=" "341×3="
Replace-Text "646×8=" "323×9="
Replace-Text "580×3=" "511×4="
Replace-Text "808×5=" "904×7="
Replace-Text "840×8=" "998×5="
Replace-Text "105×4=" "947×5="
Replace-Text "597×2=" "865×9="
Replace-Text "831×6=" "742×9="
Replace-Text "635×7=" "858×5="
Replace-Text "686×9=" "639×4="
Replace-Text "707×5=" "877×9="
Replace-Text "522×7=" "201×2="
Replace-Text "838×3=" "366×3="
Replace-Text "561×7=" "491×7="
Replace-Text "679×4=" "883×4="
Replace-Text "678×6=" "605×7="
